# TC05_Search_product_in_Catalog.xlsx — "Changes done for Kaman new UI - header & footer"
#
# Sheet "TC05_Search_product_in_Catalog" (test steps):
#   - Row 6: Object column (C6) changes from the old header-verification
#     object to the new "SearchDimensions" object.
#   - Row 7: the VERIFY_TEXT_PRESENT step is replaced by a
#     VERIFY_WEBELEMENT_PRESENT step against "ValidSearchPagination",
#     driven by data descriptor "Pagination".
#   - Row 8 (the old ValidSeachImg verification row) is removed entirely.
#
# Sheet "Testdata" (test data):
#   - B3 becomes the numeric value 200 (was the text "203kdd").
#   - A5/B5 becomes "Pagination" / TRUE (was "validSearchText" / the old
#     product description text).
#   - Two new rows are appended describing the JS element types used by
#     the new UI: EleType1/JSElement and EleType2/JSElement.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: TC05_Search_product_in_Catalog
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("TC05_Search_product_in_Catalog")

$ws1.Range("C6").Value = "SearchDimensions"

$ws1.Range("B7").Value = "VERIFY_WEBELEMENT_PRESENT"
$ws1.Range("C7").Value = "ValidSearchPagination"
$ws1.Range("E7").Value = "Pagination"

# Drop the old row 8 (ValidSeachImg check) — rows shift up, dimension becomes A1:E7
$ws1.Rows(8).Delete()

# Matches the post-edit selection captured in the workbook (rows 3-7 selected)
$ws1.Range("A3:XFD7").Select()

# ---------------------------------------------------------------------
# Sheet 2: Testdata
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Testdata")

$ws2.Range("B3").Value = 200

$ws2.Range("A5").Value = "Pagination"
$ws2.Range("B5").Value = $true

$ws2.Range("A7").Value = "EleType1"
$ws2.Range("B7").Value = "JSElement"

$ws2.Range("A8").Value = "EleType2"
$ws2.Range("B8").Value = "JSElement"

$ws2.Range("B5").Select()
